$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: insert a new "Grade" column at C, and remove the old
#     "PS 6" column (which, after the insert, has shifted from H to I). Every
#     existing formula/merge/format that lived in C:N slides over automatically.
$ws.Columns("C:C").Insert()
$ws.Columns("I:I").Delete()

# --- New column C content: a manual "Participation" / "Grade" entry column
#     that O3 (the "Participation" score) now reads from, instead of a
#     hard-coded 100.
$ws.Range("C1").Value = "Participation"
$ws.Range("C2").Value = "Grade"
$ws.Range("C5").Value = 100
$ws.Range("O3").Formula = "=C3"

# --- Formatting for the new column ---------------------------------------
# Header cell C1: bold, centered, light accent fill, no border (matches the
# other top-row group headers such as the old C1 "Problem Sets").
$c1 = $ws.Range("C1")
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108
$c1.Interior.Color = 14083579

# Sub-header C2 "Grade": bold, centered, boxed, same light accent fill.
$c2 = $ws.Range("C2")
$c2.Font.Bold = $true
$c2.HorizontalAlignment = -4108
$c2.Interior.Color = 14083579
$c2.Borders.LineStyle = 1

# Data-entry cell C3: centered, boxed, same light accent fill (keeps the
# "you type here" visual grouping of the new column).
$c3 = $ws.Range("C3")
$c3.HorizontalAlignment = -4108
$c3.Interior.Color = 14083579
$c3.Borders.LineStyle = 1

# Possible-points cell C5: italic (matches the rest of row 5), centered.
$c5 = $ws.Range("C5")
$c5.Font.Italic = $true
$c5.HorizontalAlignment = -4108

# Column C width, close to the author's 12.29-char column.
$ws.Columns("C:C").ColumnWidth = 11.5

# --- Page setup / view bits to mirror the rest of the diff -----------------
$ws.PageSetup.Orientation = 1
$ws.Range("Q3").Select()
